# NPC.xlsx ("Property1" sheet) gains a new "AIOwnerID" column, inserted
# immediately before the existing "NPCType" column (which, together with
# everything to its right, shifts one column over - NPCType: AI->AJ,
# DescID: AJ->AK, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Column 35 = AI. Insert a fresh column there; columns AI..AK (NPCType,
# DescID, ...) all shift right by one to AJ..AK.. and inherit the
# left neighbour's (AH / MasterID) formatting, same as Excel does.
$ws.Columns.Item(35).Insert() | Out-Null

# New column's width should match its left neighbour (AH / MasterID).
$ws.Columns.Item(35).ColumnWidth = $ws.Columns.Item(34).ColumnWidth

# Header row.
$ws.Range("AI1").Value = "AIOwnerID"

# Metadata rows that describe each column (Type/Public/Private/Save/Cache/
# Ref/Force/Upload) - mirror the MasterID column's (AH) metadata, since
# AIOwnerID is the same kind of "reference id" field.
$ws.Range("AI2").Value = "object"
$ws.Range("AI3").Value = 1
$ws.Range("AI4").Value = 1
$ws.Range("AI5").Value = 1
$ws.Range("AI6").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AI9").Value = 0

# Desc row.
$ws.Range("AI10").Value = "AI"

# Data rows: default every existing NPC row's new AIOwnerID value to 0.
for ($r = 11; $r -le 62; $r++) {
    $ws.Cells.Item($r, 35).Value = 0
}

# Match the author's final selection/active cell.
$ws.Range("AI11").Select() | Out-Null
